# Price-tracker update: append the 2026-02-07 scrape as a new row.
#
# The sheet stores every column (even numeric-looking price/discount
# values) as shared-string TEXT, not native numbers/dates. Plain
# `.Value = "..."` assignment lets Excel auto-coerce a "YYYY-MM-DD"
# string into a date serial and digit strings into numbers, so the
# range is first forced to Text format ("@") to keep the new cells as
# literal text, then the format is reset to the workbook's default
# "Normal" style so the new row doesn't carry a stray explicit style
# (matching the look of every other row in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 37

$rng = $ws.Range("A$newRow`:D$newRow")
$rng.NumberFormat = "@"

$ws.Range("A$newRow").Value = "2026-02-07"
$ws.Range("B$newRow").Value = "147300"
$ws.Range("C$newRow").Value = "43"
$ws.Range("D$newRow").Value = "0"

$rng.Style = "Normal"
